$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.019786254444271
$ws.Cells.Item(2, 4).Value = 1.021667955964333
$ws.Cells.Item(2, 5).Value = 1.020882804103955
$ws.Cells.Item(2, 6).Value = 1.031234167892764
$ws.Cells.Item(2, 9).Value = 1.030286179517262
$ws.Cells.Item(2, 10).Value = 1.024987050163701
$ws.Cells.Item(2, 11).Value = 1.024504607845297
$ws.Cells.Item(2, 12).Value = 1.023721774825712
$ws.Cells.Item(2, 13).Value = 1.034042875587482
$ws.Cells.Item(2, 14).Value = 1.012421919586851
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.020642677840024
$ws.Cells.Item(3, 4).Value = 1.022396566376474
$ws.Cells.Item(3, 5).Value = 1.021606293077644
$ws.Cells.Item(3, 6).Value = 1.032334060470297
$ws.Cells.Item(3, 9).Value = 1.030422133234464
$ws.Cells.Item(3, 10).Value = 1.025480865545884
$ws.Cells.Item(3, 11).Value = 1.025039933399137
$ws.Cells.Item(3, 12).Value = 1.024251825196676
$ws.Cells.Item(3, 13).Value = 1.034950508031661
$ws.Cells.Item(3, 14).Value = 1.012585179699097
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.021196995546405
$ws.Cells.Item(4, 4).Value = 1.02286847738142
$ws.Cells.Item(4, 5).Value = 1.022074981723298
$ws.Cells.Item(4, 6).Value = 1.033045795754188
$ws.Cells.Item(4, 9).Value = 1.030508115131108
$ws.Cells.Item(4, 10).Value = 1.025799927803328
$ws.Cells.Item(4, 11).Value = 1.025386110327352
$ws.Cells.Item(4, 12).Value = 1.024594679048059
$ws.Cells.Item(4, 13).Value = 1.03553725436047
$ws.Cells.Item(4, 14).Value = 1.01269064711715
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.021430066470391
$ws.Cells.Item(5, 4).Value = 1.023066975319
$ws.Cells.Item(5, 5).Value = 1.022272146883831
$ws.Cells.Item(5, 6).Value = 1.033345016166604
$ws.Cells.Item(5, 9).Value = 1.03054378493475
$ws.Cells.Item(5, 10).Value = 1.025933948175471
$ws.Cells.Item(5, 11).Value = 1.025531590620168
$ws.Cells.Item(5, 12).Value = 1.02473878394115
$ws.Cells.Item(5, 13).Value = 1.035783789538469
$ws.Cells.Item(5, 14).Value = 1.012734943882475
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.021469202165328
$ws.Cells.Item(6, 4).Value = 1.023100310200547
$ws.Cells.Item(6, 5).Value = 1.022305259270193
$ws.Cells.Item(6, 6).Value = 1.033395256978118
$ws.Cells.Item(6, 9).Value = 1.030549746058976
$ws.Cells.Item(6, 10).Value = 1.025956444105272
$ws.Cells.Item(6, 11).Value = 1.025556014290432
$ws.Cells.Item(6, 12).Value = 1.02476297794364
$ws.Cells.Item(6, 13).Value = 1.035825176061391
$ws.Cells.Item(6, 14).Value = 1.012742379045984
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.02120010971253
$ws.Cells.Item(7, 4).Value = 1.022871129303804
$ws.Cells.Item(7, 5).Value = 1.022077615751006
$ws.Cells.Item(7, 6).Value = 1.033049793925619
$ws.Cells.Item(7, 9).Value = 1.030508593627958
$ws.Cells.Item(7, 10).Value = 1.025801719037012
$ws.Cells.Item(7, 11).Value = 1.025388054449823
$ws.Cells.Item(7, 12).Value = 1.024596604707414
$ws.Cells.Item(7, 13).Value = 1.035540549098994
$ws.Cells.Item(7, 14).Value = 1.012691239177164
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.02007565401299
$ws.Cells.Item(8, 4).Value = 1.021914098908799
$ws.Cells.Item(8, 5).Value = 1.02112719726819
$ws.Cells.Item(8, 6).Value = 1.031605874924253
$ws.Cells.Item(8, 9).Value = 1.030332537147623
$ws.Cells.Item(8, 10).Value = 1.025154033931996
$ws.Cells.Item(8, 11).Value = 1.024685567355497
$ws.Cells.Item(8, 12).Value = 1.023900932861408
$ws.Cells.Item(8, 13).Value = 1.034349728215067
$ws.Cells.Item(8, 14).Value = 1.012477129661598
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.018095456588025
$ws.Cells.Item(9, 4).Value = 1.020231211871722
$ws.Cells.Item(9, 5).Value = 1.019456663159834
$ws.Cells.Item(9, 6).Value = 1.029061769260504
$ws.Cells.Item(9, 9).Value = 1.030007100834324
$ws.Cells.Item(9, 10).Value = 1.024009181471641
$ws.Cells.Item(9, 11).Value = 1.023446102345842
$ws.Cells.Item(9, 12).Value = 1.022674168984361
$ws.Cells.Item(9, 13).Value = 1.032247147806895
$ws.Cells.Item(9, 14).Value = 1.012098534214054
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.016776227962107
$ws.Cells.Item(10, 4).Value = 1.019111743834431
$ws.Cells.Item(10, 5).Value = 1.01834590335166
$ws.Cells.Item(10, 6).Value = 1.027365906701585
$ws.Cells.Item(10, 9).Value = 1.029779965716026
$ws.Cells.Item(10, 10).Value = 1.023243621491415
$ws.Cells.Item(10, 11).Value = 1.022618787586071
$ws.Cells.Item(10, 12).Value = 1.021855786914954
$ws.Cells.Item(10, 13).Value = 1.030842640191943
$ws.Cells.Item(10, 14).Value = 1.011845280015157
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.016205216124841
$ws.Cells.Item(11, 4).Value = 1.018627602824886
$ws.Cells.Item(11, 5).Value = 1.017865645521788
$ws.Cells.Item(11, 6).Value = 1.026631634169836
$ws.Cells.Item(11, 9).Value = 1.029679208626027
$ws.Cells.Item(11, 10).Value = 1.022911586635257
$ws.Cells.Item(11, 11).Value = 1.022260326549891
$ws.Cells.Item(11, 12).Value = 1.021501304241113
$ws.Cells.Item(11, 13).Value = 1.030233818847158
$ws.Cells.Item(11, 14).Value = 1.011735419286414
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.015993151505803
$ws.Cells.Item(12, 4).Value = 1.018447862300008
$ws.Cells.Item(12, 5).Value = 1.017687364302221
$ws.Cells.Item(12, 6).Value = 1.026358900055878
$ws.Cells.Item(12, 9).Value = 1.029641422137962
$ws.Cells.Item(12, 10).Value = 1.022788173598774
$ws.Cells.Item(12, 11).Value = 1.022127144935987
$ws.Cells.Item(12, 12).Value = 1.021369616969531
$ws.Cells.Item(12, 13).Value = 1.030007576755919
$ws.Cells.Item(12, 14).Value = 1.011694582421529
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.016038638502731
$ws.Cells.Item(13, 4).Value = 1.018486413111408
$ws.Cells.Item(13, 5).Value = 1.017725601321653
$ws.Cells.Item(13, 6).Value = 1.02641740210937
$ws.Cells.Item(13, 9).Value = 1.029649543785055
$ws.Cells.Item(13, 10).Value = 1.022814649741597
$ws.Cells.Item(13, 11).Value = 1.022155714330585
$ws.Cells.Item(13, 12).Value = 1.021397865071254
$ws.Cells.Item(13, 13).Value = 1.030056110913112
$ws.Cells.Item(13, 14).Value = 1.011703343406041
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.016187686084954
$ws.Cells.Item(14, 4).Value = 1.018612743552221
$ws.Cells.Item(14, 5).Value = 1.017850906530466
$ws.Cells.Item(14, 6).Value = 1.02660908973018
$ws.Cells.Item(14, 9).Value = 1.029676092539882
$ws.Cells.Item(14, 10).Value = 1.022901386911879
$ws.Cells.Item(14, 11).Value = 1.022249318387727
$ws.Cells.Item(14, 12).Value = 1.021490419263043
$ws.Cells.Item(14, 13).Value = 1.030215119616843
$ws.Cells.Item(14, 14).Value = 1.011732044302265
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.016279523838022
$ws.Cells.Item(15, 4).Value = 1.018690592008328
$ws.Cells.Item(15, 5).Value = 1.017928125550008
$ws.Cells.Item(15, 6).Value = 1.026727195789469
$ws.Cells.Item(15, 9).Value = 1.029692402311873
$ws.Cells.Item(15, 10).Value = 1.022954817906201
$ws.Cells.Item(15, 11).Value = 1.022306986566065
$ws.Cells.Item(15, 12).Value = 1.02154744278477
$ws.Cells.Item(15, 13).Value = 1.030313077057458
$ws.Cells.Item(15, 14).Value = 1.011749723946674
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.016814128928221
$ws.Cells.Item(16, 4).Value = 1.019143887360904
$ws.Cells.Item(16, 5).Value = 1.018377791535698
$ws.Cells.Item(16, 6).Value = 1.027414638971532
$ws.Cells.Item(16, 9).Value = 1.029786601996002
$ws.Cells.Item(16, 10).Value = 1.023265646208081
$ws.Cells.Item(16, 11).Value = 1.022642572746463
$ws.Cells.Item(16, 12).Value = 1.021879310372223
$ws.Cells.Item(16, 13).Value = 1.030883031787848
$ws.Cells.Item(16, 14).Value = 1.011852566926063
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.017149533035071
$ws.Cells.Item(17, 4).Value = 1.019428388165019
$ws.Cells.Item(17, 5).Value = 1.018660045859083
$ws.Cells.Item(17, 6).Value = 1.027845866763966
$ws.Cells.Item(17, 9).Value = 1.029845047265693
$ws.Cells.Item(17, 10).Value = 1.023460476080201
$ws.Cells.Item(17, 11).Value = 1.022853016731969
$ws.Cells.Item(17, 12).Value = 1.02208745135826
$ws.Cells.Item(17, 13).Value = 1.031240372885806
$ws.Cells.Item(17, 14).Value = 1.011917024313056
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.017345190004259
$ws.Cells.Item(18, 4).Value = 1.019594390082577
$ws.Cells.Item(18, 5).Value = 1.0188247482915
$ws.Cells.Item(18, 6).Value = 1.028097399167695
$ws.Cells.Item(18, 9).Value = 1.029878905205903
$ws.Cells.Item(18, 10).Value = 1.023574064656216
$ws.Cells.Item(18, 11).Value = 1.022975743050415
$ws.Cells.Item(18, 12).Value = 1.022208845032223
$ws.Cells.Item(18, 13).Value = 1.031448740164363
$ws.Cells.Item(18, 14).Value = 1.011954601897092
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.01741190757503
$ws.Cells.Item(19, 4).Value = 1.019651002140283
$ws.Cells.Item(19, 5).Value = 1.018880919097823
$ws.Cells.Item(19, 6).Value = 1.028183165975612
$ws.Cells.Item(19, 9).Value = 1.029890410483516
$ws.Cells.Item(19, 10).Value = 1.023612786508963
$ws.Cells.Item(19, 11).Value = 1.023017585772564
$ws.Cells.Item(19, 12).Value = 1.022250235164275
$ws.Cells.Item(19, 13).Value = 1.031519777196384
$ws.Cells.Item(19, 14).Value = 1.011967411594673
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.017113545124078
$ws.Cells.Item(20, 4).Value = 1.019397857972305
$ws.Cells.Item(20, 5).Value = 1.018629755570153
$ws.Cells.Item(20, 6).Value = 1.027799599632148
$ws.Cells.Item(20, 9).Value = 1.029838800653712
$ws.Cells.Item(20, 10).Value = 1.023439578090874
$ws.Cells.Item(20, 11).Value = 1.022830440361848
$ws.Cells.Item(20, 12).Value = 1.022065120975917
$ws.Cells.Item(20, 13).Value = 1.031202040157098
$ws.Cells.Item(20, 14).Value = 1.011910110641209
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.016143794322758
$ws.Cells.Item(21, 4).Value = 1.018575539862346
$ws.Cells.Item(21, 5).Value = 1.01781400427861
$ws.Cells.Item(21, 6).Value = 1.02655264229682
$ws.Cells.Item(21, 9).Value = 1.029668284548318
$ws.Cells.Item(21, 10).Value = 1.022875847191705
$ws.Cells.Item(21, 11).Value = 1.022221755234493
$ws.Cells.Item(21, 12).Value = 1.021463164810645
$ws.Cells.Item(21, 13).Value = 1.030168298238683
$ws.Cells.Item(21, 14).Value = 1.011723593420727
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.015534274250702
$ws.Cells.Item(22, 4).Value = 1.018059043055907
$ws.Cells.Item(22, 5).Value = 1.017301734194268
$ws.Cells.Item(22, 6).Value = 1.025768673682623
$ws.Cells.Item(22, 9).Value = 1.029558987046415
$ws.Cells.Item(22, 10).Value = 1.02252094188937
$ws.Cells.Item(22, 11).Value = 1.021838859035968
$ws.Cells.Item(22, 12).Value = 1.021084595432237
$ws.Cells.Item(22, 13).Value = 1.029517772635122
$ws.Cells.Item(22, 14).Value = 1.011606150984537
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.015857372972199
$ws.Cells.Item(23, 4).Value = 1.0183327973266
$ws.Cells.Item(23, 5).Value = 1.017573238511596
$ws.Cells.Item(23, 6).Value = 1.026184266102837
$ws.Cells.Item(23, 9).Value = 1.029617125346522
$ws.Cells.Item(23, 10).Value = 1.022709127688116
$ws.Cells.Item(23, 11).Value = 1.022041857369722
$ws.Cells.Item(23, 12).Value = 1.021285290929288
$ws.Cells.Item(23, 13).Value = 1.02986268253966
$ws.Cells.Item(23, 14).Value = 1.011668425597627
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.017129806445842
$ws.Cells.Item(24, 4).Value = 1.019411653074993
$ws.Cells.Item(24, 5).Value = 1.018643442237132
$ws.Cells.Item(24, 6).Value = 1.027820505743354
$ws.Cells.Item(24, 9).Value = 1.029841623946278
$ws.Cells.Item(24, 10).Value = 1.023449021154672
$ws.Cells.Item(24, 11).Value = 1.022840641720094
$ws.Cells.Item(24, 12).Value = 1.022075211149952
$ws.Cells.Item(24, 13).Value = 1.03121936126404
$ws.Cells.Item(24, 14).Value = 1.011913234691645
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.01860723112564
$ws.Cells.Item(25, 4).Value = 1.020665851932065
$ws.Cells.Item(25, 5).Value = 1.019888026408366
$ws.Cells.Item(25, 6).Value = 1.029719447113411
$ws.Cells.Item(25, 9).Value = 1.030093030956484
$ws.Cells.Item(25, 10).Value = 1.024305567504825
$ws.Cells.Item(25, 11).Value = 1.023766715685417
$ws.Cells.Item(25, 12).Value = 1.022991416781212
$ws.Cells.Item(25, 13).Value = 1.03279120926875
$ws.Cells.Item(25, 14).Value = 1.01219656287921